$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two oldest years (2008年, 2009年): those rows carried almost no
# data in the refreshed source extract. Deleting row 2 twice shifts every
# remaining row (2010年..2020年) up by two positions, turning the old
# A1:O14 block into A1:O12.
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()

# Append the newly published 2021年 figures as row 13.
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 3812
$ws.Range("C13").Value = 2945
$ws.Range("D13").Value = 1191.93
$ws.Range("E13").Value = 3682.8
$ws.Range("F13").Value = 2610
$ws.Range("G13").Value = 3476.23
$ws.Range("H13").Value = 3.7663
$ws.Range("I13").Value = 5491.27
$ws.Range("J13").Value = 143
$ws.Range("K13").Value = 4842
$ws.Range("L13").Value = 322.2617
$ws.Range("M13").Value = 6234
$ws.Range("N13").Value = 8697
$ws.Range("O13").Value = 5.5898

# Carry the bordered/centered year-label formatting used by every other
# cell in column A onto the freshly written A13 cell.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
